$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first row entry
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-14 03:34:40"

# "zh-cn" sheet: Correspond Handoff / Handback datetimes for the first row entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-14 03:34:32"
$wsZhCn.Range("K2").Value = "2016-08-14 03:35:03"

# "de-de" sheet: Correspond Handoff / Handback datetimes for the first row entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-14 03:34:40"
$wsDeDe.Range("K2").Value = "2016-08-14 03:35:14"
